# Update the cryptos worksheet with refreshed prices / 1h-volume percentages
# (scraped from coinranking.com by the GitHub Actions job), and fix the row
# order for WrappedBTC / Polkadot (rows 17-18 were swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.510.97"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "2.895.75"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").Value = "2.896.75"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.83%  "
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "3.376.69"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "60.486.13"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").Value = "2.894.81"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "424.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.669"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.28%  "
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.287"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "372.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0345"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "2.646.67"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.36%  "
$ws.Range("E51").Value = "  -0.74%  "
